$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values (regenerated sval data filtering save games)
$ws.Range("B2").Value = 0.127881588408715
$ws.Range("C2").Value = 0.002777888934908601
$ws.Range("D2").Value = 0.1575252929769615
$ws.Range("E2").Value = 0.496779210170732
$ws.Range("G2").Value = 0.7849639804913171

# Row 3 values (regenerated sval data filtering save games)
$ws.Range("B3").Value = 3.230985683306322
$ws.Range("C3").Value = 1.667794583268128
$ws.Range("D3").Value = 0.8054896365839992
$ws.Range("E3").Value = 0.496779210170732
$ws.Range("G3").Value = 6.201049113329182
